$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell values (F and G columns) for rows with revised AgTests/AgPosit figures
$ws.Cells.Item(314, 6).Value = 65185
$ws.Cells.Item(320, 6).Value = 73984
$ws.Cells.Item(320, 7).Value = 3357
$ws.Cells.Item(341, 6).Value = 283576
$ws.Cells.Item(341, 7).Value = 3608
$ws.Cells.Item(345, 6).Value = 292740
$ws.Cells.Item(355, 6).Value = 222107
$ws.Cells.Item(409, 6).Value = 708507
$ws.Cells.Item(411, 6).Value = 225427
$ws.Cells.Item(413, 6).Value = 149599
$ws.Cells.Item(416, 6).Value = 671650
$ws.Cells.Item(418, 6).Value = 202374
$ws.Cells.Item(423, 6).Value = 439862
$ws.Cells.Item(425, 6).Value = 137524
$ws.Cells.Item(428, 6).Value = 102360
$ws.Cells.Item(432, 6).Value = 118198
$ws.Cells.Item(437, 6).Value = 162020
$ws.Cells.Item(439, 6).Value = 86646
$ws.Cells.Item(447, 6).Value = 64639
$ws.Cells.Item(454, 6).Value = 50809
$ws.Cells.Item(454, 7).Value = 128
$ws.Cells.Item(458, 6).Value = 67680
$ws.Cells.Item(459, 6).Value = 57725
$ws.Cells.Item(460, 6).Value = 55772
$ws.Cells.Item(462, 6).Value = 41891
$ws.Cells.Item(463, 6).Value = 44877
$ws.Cells.Item(464, 6).Value = 69832
$ws.Cells.Item(467, 6).Value = 50355
$ws.Cells.Item(470, 6).Value = 41464
$ws.Cells.Item(471, 7).Value = 51
$ws.Cells.Item(473, 6).Value = 38946
$ws.Cells.Item(474, 6).Value = 43849
$ws.Cells.Item(475, 6).Value = 34110
$ws.Cells.Item(476, 6).Value = 35151
$ws.Cells.Item(477, 6).Value = 37160
$ws.Cells.Item(477, 7).Value = 36
$ws.Cells.Item(478, 6).Value = 51066
$ws.Cells.Item(479, 6).Value = 39752
$ws.Cells.Item(480, 6).Value = 32622
$ws.Cells.Item(481, 6).Value = 42983
$ws.Cells.Item(481, 7).Value = 37
$ws.Cells.Item(482, 6).Value = 34125
$ws.Cells.Item(483, 6).Value = 63177
$ws.Cells.Item(483, 7).Value = 35

# Append new daily rows 484-488
$ws.Cells.Item(484, 1).Value = 44378
$ws.Cells.Item(484, 2).Value = 391676
$ws.Cells.Item(484, 3).Value = 5442
$ws.Cells.Item(484, 4).Value = 17
$ws.Cells.Item(484, 5).Value = 12511
$ws.Cells.Item(484, 6).Value = 8104
$ws.Cells.Item(484, 7).Value = 11

$ws.Cells.Item(485, 1).Value = 44379
$ws.Cells.Item(485, 2).Value = 391696
$ws.Cells.Item(485, 3).Value = 8272
$ws.Cells.Item(485, 4).Value = 20
$ws.Cells.Item(485, 5).Value = 12513
$ws.Cells.Item(485, 6).Value = 12191
$ws.Cells.Item(485, 7).Value = 14

$ws.Cells.Item(486, 1).Value = 44380
$ws.Cells.Item(486, 2).Value = 391717
$ws.Cells.Item(486, 3).Value = 4675
$ws.Cells.Item(486, 4).Value = 21
$ws.Cells.Item(486, 5).Value = 12513
$ws.Cells.Item(486, 6).Value = 7544
$ws.Cells.Item(486, 7).Value = 4

$ws.Cells.Item(487, 1).Value = 44381
$ws.Cells.Item(487, 2).Value = 391720
$ws.Cells.Item(487, 3).Value = 1187
$ws.Cells.Item(487, 4).Value = 3
$ws.Cells.Item(487, 5).Value = 12513
$ws.Cells.Item(487, 6).Value = 5745
$ws.Cells.Item(487, 7).Value = 7

$ws.Cells.Item(488, 1).Value = 44382
$ws.Cells.Item(488, 2).Value = 391735
$ws.Cells.Item(488, 3).Value = 2889
$ws.Cells.Item(488, 4).Value = 15
$ws.Cells.Item(488, 5).Value = 12514
$ws.Cells.Item(488, 6).Value = 4911
$ws.Cells.Item(488, 7).Value = 5

